$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one weekly price-report row per record. A new weekly
# observation needs to be inserted right after the existing row for this
# market/date cluster (row 116), pushing the remaining rows (116-180)
# down by one (to 117-181) - exactly like Excel's native "Insert Row".
$ws.Rows("116").Insert()

# Populate the newly inserted row 116 with the new weekly record.
$ws.Cells.Item(116, 1).Value = 4
$ws.Cells.Item(116, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(116, 3).Value = "Los Lagos"
$ws.Cells.Item(116, 4).Value = 44719
$ws.Cells.Item(116, 5).Value = 10
$ws.Cells.Item(116, 6).Value = 100112009
$ws.Cells.Item(116, 7).Value = "Acelga"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 90
$ws.Cells.Item(116, 11).Value = 12000
$ws.Cells.Item(116, 12).Value = 12000
$ws.Cells.Item(116, 13).Value = 12000
$ws.Cells.Item(116, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(116, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(116, 16).Value = 1000
$ws.Cells.Item(116, 17).Value = 12
$ws.Cells.Item(116, 18).Value = "Hortaliza"
